# ***v1.12.2.4    2020/10/27  (Mark Chou) 功能修改 加入二期MTL,MTS,HID設定
# Update CSTTranSchedule matrix: rename stations, drop the 6th
# row/column (20616) and rewrite the 5x5 transfer matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear column F and row 6 entirely - table shrinks from 6x6 to 5x5
$ws.Range("A6:F6").Clear()
$ws.Range("F1:F5").Clear()

# Header row (B1:E1) - station IDs across the top
$ws.Range("B1").Value = 20411
$ws.Range("C1").Value = 20416
$ws.Range("D1").Value = 20316
$ws.Range("E1").Value = 20611

# Row labels (A2:A5) - station IDs down the side
$ws.Range("A2").Value = 20411
$ws.Range("A3").Value = 20416
$ws.Range("A4").Value = 20316
$ws.Range("A5").Value = 20611

# Row 2 (20411) values
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1

# Row 3 (20416) values
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 1

# Row 4 (20316) values
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

# Row 5 (20611) values
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Update the active selection to match the committed state
$ws.Range("H4:H5").Select()
